$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new items (KENACOMB, TERRAMYCIN) were added to the shortage report.
# The report lists items alphabetically, so they land between "BRONCHOPRO..."
# (row 9) and "TRIACTIN..." (row 10) while the two rows that previously held
# TRIACTIN/ZURCAL/syringes shift down. Net effect on the sheet: the totals
# row and the footer row move down by two rows, and two new data rows are
# created just above the totals row, matching the same look as the other
# data rows (row 12, the last data row before the insert).

# 1) Insert two blank rows right before the totals row (row 13).
$ws.Rows("13:14").Insert()

# 2) Give the two new rows the same formatting / merged-cell layout as the
#    other data rows by copying the last data row (row 12) into them.
$ws.Range("A12:Q12").Copy($ws.Range("A13:Q13"))
$ws.Range("A12:Q12").Copy($ws.Range("A14:Q14"))

# 3) Restore the expected row heights for the two new rows.
$ws.Rows("13").RowHeight = 24.75
$ws.Rows("14").RowHeight = 25.5

# 4) Rewrite the data rows in alphabetical order, inserting the two new
#    medicines in their correct position.
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "KENACOMB TOPICAL CREAM 15 GM"
$ws.Range("H10").Value = "2:0"
$ws.Range("L10").Value = "1"
$ws.Range("N10").Value = "36.00"
$ws.Range("P10").Value = "36.0000"
$ws.Range("Q10").Value = "1:0"

$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "TERRAMYCIN EYE OINT. 5 GM"
$ws.Range("H11").Value = "1:0"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "28.00"
$ws.Range("P11").Value = "28.0000"
$ws.Range("Q11").Value = "1:0"

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "TRIACTIN 4MG 20 TAB"
$ws.Range("H12").Value = "1:0"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "46.00"
$ws.Range("P12").Value = "23.0000"
$ws.Range("Q12").Value = "0:1"

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "ZURCAL 40MG 14 GASTRO RESISTANT TAB"
$ws.Range("H13").Value = "4:0"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "96.00"
$ws.Range("P13").Value = "96.0000"
$ws.Range("Q13").Value = "1:0"

$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "سرنجات 3 سم"
$ws.Range("H14").Value = "0:0"
$ws.Range("L14").Value = "0"
$ws.Range("N14").Value = "2.00"
$ws.Range("P14").Value = "2.0000"
$ws.Range("Q14").Value = "1:0"

# 5) Update the grand total (sum of the "sell price" column) now that two
#    more items are included.
$ws.Range("P15").Value = 266.48

# 6) Update the "generated at" timestamp in the footer.
$ws.Range("A16").Value = "Thursday, 12 June, 2025 10:42 AM"
